$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44216
$ws.Range("P2").Value2 = 11545
$ws.Range("S2").Value2 = 825

# Row 3
$ws.Range("D3").Value2 = 44210
$ws.Range("M3").Value2 = 70
$ws.Range("N3").Value2 = 10000
$ws.Range("O3").Value2 = 11000
$ws.Range("P3").Value2 = 10357
$ws.Range("S3").Value2 = 740

# Row 4
$ws.Range("D4").Value2 = 44181
$ws.Range("M4").Value2 = 65
$ws.Range("N4").Value2 = 9000
$ws.Range("O4").Value2 = 10000
$ws.Range("P4").Value2 = 9462
$ws.Range("S4").Value2 = 676

# Row 6
$ws.Range("D6").Value2 = 44253
$ws.Range("M6").Value2 = 90
$ws.Range("N6").Value2 = 12000
$ws.Range("O6").Value2 = 13000
$ws.Range("P6").Value2 = 12667
$ws.Range("S6").Value2 = 905

# Row 7
$ws.Range("D7").Value2 = 44232
$ws.Range("M7").Value2 = 60
$ws.Range("N7").Value2 = 11000
$ws.Range("O7").Value2 = 12000
$ws.Range("P7").Value2 = 11583
$ws.Range("S7").Value2 = 827

# Row 8
$ws.Range("D8").Value2 = 44229
$ws.Range("M8").Value2 = 55
$ws.Range("N8").Value2 = 11000
$ws.Range("O8").Value2 = 12000
$ws.Range("P8").Value2 = 11364
$ws.Range("S8").Value2 = 812
